# Integration Testing participants and methodology
# Applies the 7 localized changes described by the XML diff:
#   1. 4.1 Unit Testing / Definition: wrap "a" in proofErr gramStart/gramEnd
#   2. 4.2 System and Integration Testing / Definition: split off trailing
#      "functionality." into its own run
#   3. 4.2 .../ Participants: fill in the previously-empty paragraph with
#      three participant names
#   4. 4.2 .../ Methodology: fill in the previously-empty paragraph with
#      the methodology text
#   5. 4.4 User Acceptance Testing / Participants: gains a lastRenderedPageBreak
#   6. 4.5 Batch Testing heading loses its lastRenderedPageBreak
#   7. 5.0 Test Schedule heading gains a lastRenderedPageBreak

$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function New-PackageXml($bodyFragment) {
    return @"
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>
<pkg:xmlData>
<w:document $wNs>
<w:body>
$bodyFragment
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
}

# ---------------------------------------------------------------------------
# 1. 4.1 Unit Testing -- Definition paragraph: wrap "a" with proofErr markers
# ---------------------------------------------------------------------------
$rng = $d.Content
$ok = $rng.Find.Execute(
    "Testing individual units of the product to get a accurate test of the products specific component",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "Could not find unit-testing definition paragraph" }
$para = $rng.Paragraphs(1).Range

$frag1 = '<w:p><w:r><w:t xml:space="preserve">Testing individual units of the product to get </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>a</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> accurate test of the products specific component</w:t></w:r></w:p>'
$para.InsertXML((New-PackageXml $frag1))

# ---------------------------------------------------------------------------
# 2. 4.2 System and Integration Testing -- Definition paragraph: split off
#    "functionality." into its own trailing run
# ---------------------------------------------------------------------------
$rng = $d.Content
$ok = $rng.Find.Execute(
    "Once unit testing is completed it will move on to the integration testing where all the components tested in unit testing will be tested as a whole and see how each component works with each other to see if there is an conflicting code ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "Could not find integration-testing definition paragraph" }
$para = $rng.Paragraphs(1).Range

$frag2 = '<w:p><w:r><w:t xml:space="preserve">Once unit testing is completed it will move on to the integration testing where all the components tested in unit testing will be tested as a whole and see how each component works with each other to see if there is an conflicting </w:t></w:r><w:r><w:t>functionality.</w:t></w:r></w:p>'
$para.InsertXML((New-PackageXml $frag2))

# ---------------------------------------------------------------------------
# 3. 4.2 .../ Participants: fill the empty paragraph with three names
# ---------------------------------------------------------------------------
$rng = $d.Content
$ok = $rng.Find.Execute("4.2 System and Integration Testing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "Could not find 4.2 heading" }
$p = $rng.Paragraphs(1)
$p = $p.Next()  # Definition:
$p = $p.Next()  # definition text
$p = $p.Next()  # Participants:
$participantsEmpty = $p.Next()  # empty paragraph to fill in (step 3)

$frag3 = '<w:p><w:r><w:t xml:space="preserve">May </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Mcgee</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ellisha</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Osborne</w:t></w:r></w:p><w:p><w:r><w:t>Yusra Cross</w:t></w:r></w:p>'
$participantsEmpty.Range.InsertXML((New-PackageXml $frag3))

# ---------------------------------------------------------------------------
# 4. 4.2 .../ Methodology: fill the empty paragraph with the methodology text
#    (re-find fresh -- the paragraph objects captured above are stale now
#    that step 3 inserted new paragraphs into the document)
# ---------------------------------------------------------------------------
$rng = $d.Content
$ok = $rng.Find.Execute("4.2 System and Integration Testing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "Could not find 4.2 heading (pass 2)" }
$p = $rng.Paragraphs(1)
$p = $p.Next()  # Definition:
$p = $p.Next()  # definition text
$p = $p.Next()  # Participants:
$p = $p.Next()  # May Mcgee
$p = $p.Next()  # Ellisha Osborne
$p = $p.Next()  # Yusra Cross
$p = $p.Next()  # Methodology:
$methodologyEmpty = $p.Next()  # empty paragraph to fill in (step 4)

$frag4 = '<w:p><w:r><w:t xml:space="preserve">With the knowledge from unit testing the participants will test the game </w:t></w:r><w:r><w:t>with</w:t></w:r><w:r><w:t xml:space="preserve"> their previous knowledge they will notice any loss or errors in functionality. We will see how the enemies and players interact with the background and platforms. We will test</w:t></w:r><w:r><w:t xml:space="preserve"> how the player and enemy/boss interact with one another and vice versa. We test that both the save and delete game work with one another.</w:t></w:r></w:p>'
$methodologyEmpty.Range.InsertXML((New-PackageXml $frag4))

# ---------------------------------------------------------------------------
# 5. 4.4 User Acceptance Testing -- Participants: gains lastRenderedPageBreak
# ---------------------------------------------------------------------------
$rng = $d.Content
$ok = $rng.Find.Execute("4.4 User Acceptance Testing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "Could not find 4.4 heading" }
$p = $rng.Paragraphs(1)
$p = $p.Next()  # Definition:
$p = $p.Next()  # empty
$p = $p.Next()  # Participants:
$para = $p.Range

$frag5 = '<w:p w14:paraId="50AA857D" w14:textId="77777777" w:rsidR="009A74C6" w:rsidRDefault="009A74C6" w:rsidP="009A74C6" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Participants:</w:t></w:r></w:p>'
$para.InsertXML((New-PackageXml $frag5))

# ---------------------------------------------------------------------------
# 6. 4.5 Batch Testing heading loses its lastRenderedPageBreak
# ---------------------------------------------------------------------------
$rng = $d.Content
$ok = $rng.Find.Execute("4.5 Batch Testing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "Could not find 4.5 heading" }
$para = $rng.Paragraphs(1).Range

$frag6 = '<w:p w14:paraId="57EA709F" w14:textId="3A4A93E6" w:rsidR="002F25A6" w:rsidRDefault="002F25A6" w:rsidP="00081CE7" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="Heading2"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:bookmarkStart w:id="12" w:name="_Toc40357032"/><w:r w:rsidRPr="00A80156"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>4.5 Batch Testing</w:t></w:r><w:bookmarkEnd w:id="12"/></w:p>'
$para.InsertXML((New-PackageXml $frag6))

# ---------------------------------------------------------------------------
# 7. 5.0 Test Schedule heading gains a lastRenderedPageBreak
# ---------------------------------------------------------------------------
$rng = $d.Content
$ok = $rng.Find.Execute("5.0 Test Schedule", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "Could not find 5.0 heading" }
$para = $rng.Paragraphs(1).Range

$frag7 = '<w:p w14:paraId="0E651FB6" w14:textId="6DC68B8A" w:rsidR="002F25A6" w:rsidRPr="00A80156" w:rsidRDefault="002F25A6" w:rsidP="00081CE7" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="Heading1"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:bookmarkStart w:id="15" w:name="_Toc40357035"/><w:r w:rsidRPr="00A80156"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:lastRenderedPageBreak/><w:t>5.0 Test Schedule</w:t></w:r><w:bookmarkEnd w:id="15"/></w:p>'
$para.InsertXML((New-PackageXml $frag7))

Write-Host "Done applying all edits."
